$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A24").Value = "distributiveProperty"
$ws.Range("B24").Value = "Distributive Property"

$ws.Range("A25").Value = "areaOfRectangle"
$ws.Range("B25").Value = "Area of Rectangle"
$ws.Range("B25").WrapText = $true

$ws.Range("A24").Select()
